# Updated keywords and implemented Screenshot attachments in reporting
#
# Inserts a new "GooglePage" worksheet (with a searchBox/css/#APjFqb locator
# row) between DashboardPage and ExamplePage, and adjusts the active
# sheet/selection state to match.

$wb = $excel.ActiveWorkbook

# --- Add the new GooglePage sheet right after DashboardPage ---------------
$dashboard = $wb.Worksheets.Item("DashboardPage")
$googlePage = $wb.Worksheets.Add($null, $dashboard)
$googlePage.Name = "GooglePage"

# Header row (same locator table layout as the other sheets)
$googlePage.Range("A1").Value = "LocatorName"
$googlePage.Range("B1").Value = "LocatorType"
$googlePage.Range("C1").Value = "LocatorValue"

# New locator: the Google search box
$googlePage.Range("A2").Value = "searchBox"
$googlePage.Range("B2").Value = "css"
$googlePage.Range("C2").Value = "#APjFqb"

# --- Update selections / active sheet --------------------------------------
# DashboardPage: selection now spans A1:D2
$dashboard.Range("A1:D2").Select()

# GooglePage becomes the active tab, with C2 as the active/selected cell
$googlePage.Select()
$googlePage.Range("C2").Select()
